$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (hypothesis) text updates
$ws.Range("B2").Value = "<pare>"
$ws.Range("B6").Value = "<para>"
$ws.Range("B11").Value = "<him>"
$ws.Range("B12").Value = "<him>"
$ws.Range("B13").Value = "<him>"
$ws.Range("B14").Value = "<make>"
$ws.Range("B15").Value = "<out>"
$ws.Range("B16").Value = "<down>"

# Column C (score) updates
$ws.Range("C2").Value = 38
$ws.Range("C3").Value = 35
$ws.Range("C4").Value = 40
$ws.Range("C5").Value = 39
$ws.Range("C7").Value = 37
$ws.Range("C8").Value = 36
$ws.Range("C9").Value = 31
$ws.Range("C10").Value = 44
$ws.Range("C11").Value = 35
$ws.Range("C12").Value = 32
$ws.Range("C13").Value = 41
$ws.Range("C14").Value = 30
$ws.Range("C15").Value = 36
